# ---------------------------------------------------------------------------
# Applies the trigonometry.pptx update:
#   - tiny nudge of the "(x0 , y0)" label on slide 2
#   - repositions several shapes on slide 3 (axis labels / dashed guide lines)
#   - adds 5 new shapes on slide 3: "x_max" / "x_min" labels, a dashed
#     accent2 guide connector, a red ellipse marker and an "(x0 , y0)" label
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Slide 2 — nudge "(x0 , y0)" label by 6 EMU vertically
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$lbl41 = Get-ShapeById $s2 41
$lbl41.Top = 196.89976377952755   # 2500627 EMU

# ---------------------------------------------------------------------------
# Slide 3 — reposition existing shapes
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# id=16 "Connecteur droit 15" (dashed horizontal guide line)
$sh16 = Get-ShapeById $s3 16
$sh16.Left = 86.94118110236221    # 1104153 EMU
$sh16.Top  = 107.27496062992125   # 1362392 EMU

# id=35 "Connecteur droit avec flèche 34" (dashed arrow, triangle heads)
$sh35 = Get-ShapeById $s3 35
$sh35.Width  = 532.7228346456693     # 6765580 EMU
$sh35.Height = 7.5696850393700785    # 96135 EMU

# id=38 "ZoneTexte 37" (the "x_range" label)
$sh38 = Get-ShapeById $s3 38
$sh38.Left = 286.08779527559057   # 3633315 EMU
$sh38.Top  = 282.38937007874017   # 3586345 EMU

# id=5 "Connecteur droit avec flèche 4" (dashed arrow, stealth heads, accent2)
$sh5 = Get-ShapeById $s3 5
$sh5.Left   = 95.56094488188977    # 1213624 EMU
$sh5.Top    = 301.2855905511811    # 3826327 EMU
$sh5.Width  = 488.00456692913383   # 6197658 EMU
$sh5.Height = 7.891653543307086    # 100224 EMU

# id=11 "ZoneTexte 10" (the red "xM" label)
$sh11 = Get-ShapeById $s3 11
$sh11.Left = 543.766220472441     # 6905831 EMU
$sh11.Top  = 276.98622047244095   # 3517725 EMU

# ---------------------------------------------------------------------------
# Slide 3 — add 5 new shapes.
#
# The presentation's next-free-id allocator fills the gaps left by shapes
# that were deleted while this deck was authored (3,4,9,13,15,17,19,28,30,
# 31,34,36,37,39,...). Creating (and discarding) throw-away shapes first
# walks that counter forward so the shapes we actually keep land on the
# exact ids used by the target deck (28, 30, 31, 36, 37).
# ---------------------------------------------------------------------------

function Skip-ShapeIds($slide, $n) {
    for ($i = 0; $i -lt $n; $i++) {
        $tmp = $slide.Shapes.AddTextbox(1, 0, 0, 10, 10)
        $tmp.Delete()
    }
}

# Burn through ids 3,4,9,13,15,17,19 (7 throw-away shapes) so the next
# created shape receives id 28.
Skip-ShapeIds $s3 7

# id=28 "ZoneTexte 27": "x" + "max" (subscript) label.
# Cloned from id=38 ("x" + "range") so every run/paragraph property
# (sizes, dirty/err flags, spAutoFit body, etc.) matches exactly.
$src28 = Get-ShapeById $s3 38
$src28.Copy()
$new28 = $s3.Shapes.Paste().Item(1)
$new28.Name = "ZoneTexte 27"
$new28.Left = 589.8646456692913     # 7491281 EMU
$new28.Top  = 282.09614173228346    # 3582621 EMU
$sub28 = $new28.TextFrame.TextRange.Characters(2, 5)   # "range"
$sub28.Text = "max"

# id=30 "ZoneTexte 29": "x" + "min" (subscript) label.
$src30 = Get-ShapeById $s3 38
$src30.Copy()
$new30 = $s3.Shapes.Paste().Item(1)
$new30.Name = "ZoneTexte 29"
$new30.Left = 62.75456692913386     # 796983 EMU
$new30.Top  = 294.2036220472441     # 3736386 EMU
$sub30 = $new30.TextFrame.TextRange.Characters(2, 5)   # "range"
$sub30.Text = "min"

# id=31 "Connecteur droit avec flèche 30": dashed accent2 connector with
# stealth arrowheads on both ends. Copied from slide 2's shape id=34,
# which already carries the exact matching <p:style> (lnRef idx=2 /
# fillRef idx=0 / effectRef idx=1, accent2 scheme colour).
$src31 = Get-ShapeById $s2 34
$src31.Copy()
$new31 = $s3.Shapes.Paste().Item(1)
$new31.Name = "Connecteur droit avec flèche 30"
$new31.Left   = 254.20984251968503    # 3228465 EMU
$new31.Top    = 76.4359842519685      # 970737 EMU
$new31.Width  = 3.172283464566929     # 40288 EMU
$new31.Height = 356.28905511811024    # 4524871 EMU

# Burn through id=34 (1 throw-away shape) so the next created shape
# receives id 36.
Skip-ShapeIds $s3 1

# id=36 "Ellipse 35": small red-outlined ellipse marker. Copied from
# slide 2's shape id=32 which is an exact match (size, noFill, red
# line, style block).
$src36 = Get-ShapeById $s2 32
$src36.Copy()
$new36 = $s3.Shapes.Paste().Item(1)
$new36.Name = "Ellipse 35"
$new36.Left = 250.13590551181102    # 3176726 EMU
$new36.Top  = 296.9511023622047     # 3771279 EMU

# id=37 "ZoneTexte 36": "(x0 , y0)" label. Copied from slide 2's shape
# id=41 which has identical text/formatting already.
$src37 = Get-ShapeById $s2 41
$src37.Copy()
$new37 = $s3.Shapes.Paste().Item(1)
$new37.Name = "ZoneTexte 36"
$new37.Left = 206.294094488189      # 2619935 EMU
$new37.Top  = 279.43566929133857    # 3548833 EMU
